$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 100.85714
$ws.Range("I11").Value = 100.85714
$ws.Range("K11").Value = 100.85714
$ws.Range("M11").Value = 39.14286
$ws.Range("H16").Value = 7249.75
$ws.Range("I16").Value = 7000
$ws.Range("J16").Value = 7333
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7333
$ws.Range("M16").Value = -6770
$ws.Range("N16").Value = -7793
$ws.Range("H33").Value = 21809.4
$ws.Range("J33").Value = 405.5
$ws.Range("L33").Value = 405.5
$ws.Range("N33").Value = -863.5
$ws.Range("H100").Value = 1184.4375
$ws.Range("I100").Value = 1182.2142
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 1182.2142
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -641.2141999999999
$ws.Range("N100").Value = -2282
$ws.Range("H107").Value = 704.1667
$ws.Range("I107").Value = 704.1667
$ws.Range("K107").Value = 704.1667
$ws.Range("M107").Value = 1215.8333
$ws.Range("H132").Value = 4777.6294
$ws.Range("I132").Value = 4918.885
$ws.Range("K132").Value = 14756.655
$ws.Range("M132").Value = -12226.655
$ws.Range("H138").Value = 1737.6666
$ws.Range("J138").Value = 2886.818
$ws.Range("L138").Value = 8660.454000000002
$ws.Range("N138").Value = -18940.454

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 19000
$ws.Range("J37").Value = 19000
$ws.Range("L37").Value = 19000
$ws.Range("N37").Value = -19546
$ws.Range("H44").Value = 29900
$ws.Range("J44").Value = 29900
$ws.Range("L44").Value = 29900
$ws.Range("N44").Value = -30876
$ws.Range("H61").Value = 3339567.5
$ws.Range("I61").Value = 4766947
$ws.Range("K61").Value = 4766947
$ws.Range("M61").Value = -4766735
$ws.Range("H74").Value = 4694.2964
$ws.Range("I74").Value = 2659.1875
$ws.Range("J74").Value = 7654.4546
$ws.Range("K74").Value = 2659.1875
$ws.Range("L74").Value = 7654.4546
$ws.Range("M74").Value = -1785.1875
$ws.Range("N74").Value = -9402.454600000001
$ws.Range("H77").Value = 4694.2964
$ws.Range("I77").Value = 2659.1875
$ws.Range("J77").Value = 7654.4546
$ws.Range("K77").Value = 13295.9375
$ws.Range("L77").Value = 38272.273
$ws.Range("M77").Value = -8927.9375
$ws.Range("N77").Value = -47008.273
$ws.Range("H102").Value = 4766.524
$ws.Range("I102").Value = 3818.9375
$ws.Range("K102").Value = 3818.9375
$ws.Range("M102").Value = -2196.9375
$ws.Range("H132").Value = 3911.762
$ws.Range("I132").Value = 3224.1667
$ws.Range("K132").Value = 9672.500100000001
$ws.Range("M132").Value = -7142.500100000001
$ws.Range("H135").Value = 171875
$ws.Range("J135").Value = 171875
$ws.Range("L135").Value = 171875
$ws.Range("N135").Value = -182015
$ws.Range("H136").Value = 3339567.5
$ws.Range("I136").Value = 4766947
$ws.Range("K136").Value = 14300841
$ws.Range("M136").Value = -14298291

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1392
$ws.Range("I7").Value = 888
$ws.Range("K7").Value = 888
$ws.Range("M7").Value = -775
$ws.Range("H20").Value = 2849.8235
$ws.Range("I20").Value = 3024.4
$ws.Range("J20").Value = 2600.4285
$ws.Range("K20").Value = 3024.4
$ws.Range("L20").Value = 2600.4285
$ws.Range("M20").Value = -2777.4
$ws.Range("N20").Value = -3094.4285
$ws.Range("H86").Value = 10527297
$ws.Range("I86").Value = 1103.2307
$ws.Range("K86").Value = 1103.2307
$ws.Range("M86").Value = 19.76929999999993
$ws.Range("H89").Value = 10527297
$ws.Range("I89").Value = 1103.2307
$ws.Range("K89").Value = 5516.1535
$ws.Range("M89").Value = 99.84649999999965

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4202.5
$ws.Range("I3").Value = 4202.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4202.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4089.5
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 3643.9167
$ws.Range("I22").Value = 4692.778
$ws.Range("K22").Value = 4692.778
$ws.Range("M22").Value = -4342.778
$ws.Range("H31").Value = 11097
$ws.Range("I31").Value = 9755.4
$ws.Range("K31").Value = 9755.4
$ws.Range("M31").Value = -9460.4
$ws.Range("H34").Value = 11097
$ws.Range("I34").Value = 9755.4
$ws.Range("K34").Value = 9755.4
$ws.Range("M34").Value = -9553.4
$ws.Range("H58").Value = 7918.091
$ws.Range("I58").Value = 3420
$ws.Range("J58").Value = 11666.5
$ws.Range("K58").Value = 3420
$ws.Range("L58").Value = 11666.5
$ws.Range("M58").Value = -3217
$ws.Range("N58").Value = -12072.5
$ws.Range("H105").Value = 1136.6957
$ws.Range("I105").Value = 1177.2
$ws.Range("J105").Value = 866.6667
$ws.Range("K105").Value = 1177.2
$ws.Range("L105").Value = 866.6667
$ws.Range("M105").Value = 569.8
$ws.Range("N105").Value = -4360.6667
$ws.Range("H132").Value = 2037.25
$ws.Range("I132").Value = 999.5
$ws.Range("K132").Value = 2998.5
$ws.Range("M132").Value = -468.5
$ws.Range("H136").Value = 7918.091
$ws.Range("I136").Value = 3420
$ws.Range("J136").Value = 11666.5
$ws.Range("K136").Value = 10260
$ws.Range("L136").Value = 34999.5
$ws.Range("M136").Value = -7710
$ws.Range("N136").Value = -40099.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9936237
$ws.Range("I4").Value = 5056813
$ws.Range("K4").Value = 15170439
$ws.Range("M4").Value = -15170327
$ws.Range("H6").Value = 2189.7273
$ws.Range("I6").Value = 493.1111
$ws.Range("J6").Value = 9824.5
$ws.Range("K6").Value = 1479.3333
$ws.Range("L6").Value = 29473.5
$ws.Range("M6").Value = -1366.3333
$ws.Range("N6").Value = -29699.5
$ws.Range("H47").Value = 2133
$ws.Range("H97").Value = 1390.5652
$ws.Range("I97").Value = 1364.4286
$ws.Range("J97").Value = 1431.2222
$ws.Range("K97").Value = 4093.2858
$ws.Range("L97").Value = 4293.6666
$ws.Range("M97").Value = -3597.2858
$ws.Range("N97").Value = -5285.6666
$ws.Range("H98").Value = 345
$ws.Range("I98").Value = 345
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1035
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 463
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 2861.25
$ws.Range("J122").Value = 2861.25
$ws.Range("L122").Value = 25751.25
$ws.Range("N122").Value = -30651.25
$ws.Range("H134").Value = 1924
$ws.Range("I134").Value = 1924
$ws.Range("K134").Value = 5772
$ws.Range("M134").Value = -702

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 214804.2
$ws.Range("J24").Value = 16449.111
$ws.Range("L24").Value = 16449.111
$ws.Range("N24").Value = -16795.111
$ws.Range("H33").Value = 514999.5
$ws.Range("J33").Value = 29999
$ws.Range("L33").Value = 29999
$ws.Range("N33").Value = -30503
$ws.Range("H93").Value = 37285.77
$ws.Range("J93").Value = 37285.77
$ws.Range("L93").Value = 37285.77
$ws.Range("N93").Value = -41029.77
$ws.Range("H102").Value = 2585.2646
$ws.Range("I102").Value = 2034.9615
$ws.Range("K102").Value = 2034.9615
$ws.Range("M102").Value = -412.9614999999999
$ws.Range("H132").Value = 3706.9412
$ws.Range("I132").Value = 2386.077
$ws.Range("J132").Value = 7999.75
$ws.Range("K132").Value = 7158.231000000001
$ws.Range("L132").Value = 23999.25
$ws.Range("M132").Value = -4628.231000000001
$ws.Range("N132").Value = -29059.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2917.25
$ws.Range("I7").Value = 2917.25
$ws.Range("K7").Value = 2917.25
$ws.Range("M7").Value = -2805.25
$ws.Range("H100").Value = 4171207.8
$ws.Range("I100").Value = 10003898
$ws.Range("K100").Value = 10003898
$ws.Range("M100").Value = -10003357
$ws.Range("H103").Value = 55333
$ws.Range("J103").Value = 55333
$ws.Range("L103").Value = 55333
$ws.Range("N103").Value = -57677
$ws.Range("H126").Value = 2917.25
$ws.Range("I126").Value = 2917.25
$ws.Range("K126").Value = 8751.75
$ws.Range("M126").Value = -6281.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 29779
$ws.Range("J112").Value = 29779
$ws.Range("L112").Value = 29779
$ws.Range("N112").Value = -32733
$ws.Range("H132").Value = 4676.0225
$ws.Range("I132").Value = 3897.9744
$ws.Range("K132").Value = 11693.9232
$ws.Range("M132").Value = -9163.923200000001
